# DPLKAKT071-001 - Setup Periode Bulanan - Hapus Data
# Update Regresi Tanggal 31/03/2023
# - TGL_AWAL (O2): 15/04/2023 -> 15/04/2024
# - PERIODE_BULANAN (Q2): 202305 -> 202405
# - VERIFIKASI (T2): 202305 -> 202405
# - Move the active selection from O2 to X2 (and scroll the view over)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = "15/04/2024"
$ws.Range("Q2").Value = "202405"
$ws.Range("T2").Value = 202405

$ws.Range("X2").Select()
